$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A34").Value = "2025/12/04 02:00"
$ws.Range("B34").Value = "-"
$ws.Range("C34").Value = "-"
$ws.Range("D34").Value = "-"
$ws.Range("E34").Value = "-"
$ws.Range("F34").Value = "-"
$ws.Range("G34").Value = "-"
